$d = $word.ActiveDocument

$newUrl = "https://github.com/stiven-skyward/DevOpsTraining/tree/main/WEB%20STACK%20IMPLEMENTATION%20IN%20AWS/Web%20Solution%20With%20WordPress"
$newUrlLen = $newUrl.Length

# Step 1: insert the new URL text at the very start of the document
# (position 0), ahead of the old hyperlink. Inserting at the absolute
# start of the story yields a clean run with no inherited character
# formatting (no rPr). Note: this Range object's Start/End do not
# auto-grow after InsertBefore, so compute offsets explicitly.
$insertPoint = $d.Range(0, 0)
$insertPoint.InsertBefore($newUrl)

# Step 2: remove the old hyperlink's visible text (the
# "Client-Server Architecture with MySQL" link), leaving the first
# paragraph holding only the freshly inserted, unformatted run.
$p1 = $d.Paragraphs(1)
$oldLinkRange = $d.Range($newUrlLen, $p1.Range.End - 1)
$oldLinkRange.Delete()

# Step 3: merge the now-empty-of-link first paragraph into the second
# paragraph by deleting the paragraph mark between them. The merged
# paragraph keeps the second paragraph's pPr/rPr (spacing + paragraph
# mark run formatting), while the plain new-URL run stays unformatted
# and the original run (holding the line break) is preserved intact.
$p1 = $d.Paragraphs(1)
$markRange = $d.Range($p1.Range.End - 1, $p1.Range.End)
$markRange.Delete()
